# Added hint system for area evaluate.
# New rows 85-94: hint keys/messages for "too few/many zeroes" and for each
# multiplication trick (multiples of 2 through 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: multiplication-hint rows (87-94) -------------------------------
# Column B (the hint message) is entered first for this block, then column A
# (the hint key), matching the order the strings were originally authored in.
$multMessages = @(
    "Incorrect Multiplication: Remember that in multiples of 2, simply double the number.",
    "Incorrect Multiplication: The trick with multiples of 3 is to double the number, and then add the original number.",
    "Incorrect Multiplication: In multiples of 4, double the number, and then double it again.",
    "Incorrect Multiplication: If you are having trouble with multiples of 5, try multiplying the number by 10, and then half it.",
    "Incorrect Multiplication: A good way to solve multiples of 6 is to multiply the number by 5, and then add the original number.",
    "Incorrect Multiplication: A good way to solve multiples of 7 is to multiply the number by 5, and then add the original number twice.",
    "Incorrect Multiplication: If you are having trouble with multiples of 8, try multiplying the number by 2 three times.",
    "Incorrect Multiplication: One way to do multiples of 9 is to multiply the number by 10, and then subtract it by the original number."
)
for ($i = 0; $i -lt $multMessages.Count; $i++) {
    $ws.Cells.Item(87 + $i, 2).Value = $multMessages[$i]
}

$multKeys = @("hint_mult_2", "hint_mult_3", "hint_mult_4", "hint_mult_5", "hint_mult_6", "hint_mult_7", "hint_mult_8", "hint_mult_9")
for ($i = 0; $i -lt $multKeys.Count; $i++) {
    $ws.Cells.Item(87 + $i, 1).Value = $multKeys[$i]
}

# --- Step 2: zero-count hint rows (85-86) -----------------------------------
# Column A (key) first, then column B (message) for this block.
$ws.Cells.Item(85, 1).Value = "hint_zeroes_few"
$ws.Cells.Item(86, 1).Value = "hint_zeroes_many"

$ws.Cells.Item(85, 2).Value = "Incorrect number of zeroes at the end! Looks like you missed a few zeroes."
$ws.Cells.Item(86, 2).Value = "Incorrect number of zeroes at the end! Looks like you put in too many zeroes."

# --- Step 3: MaxChars column (C) for every new row --------------------------
$ws.Cells.Item(85, 3).Value = 6
$ws.Cells.Item(86, 3).Value = 6
$ws.Cells.Item(87, 3).Value = 6
$ws.Cells.Item(88, 3).Value = 8
$ws.Cells.Item(89, 3).Value = 6
$ws.Cells.Item(90, 3).Value = 8
$ws.Cells.Item(91, 3).Value = 8
$ws.Cells.Item(92, 3).Value = 8
$ws.Cells.Item(93, 3).Value = 8
$ws.Cells.Item(94, 3).Value = 8

# --- Step 4: formatting -------------------------------------------------------
# Row 87's hint message gets the same vertically-centered style already used
# elsewhere in the sheet for longer entries.
$ws.Range("B87").VerticalAlignment = -4108

# --- Step 5: view state -------------------------------------------------------
# Reflect the newly-added area in the window scroll position/selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$ws.Range("B94").Select()
